# JPM_UPCOMING_INDEX_EVENTS.xlsx update:
# Remove the "NRP SJ / Nepi Rockcastle" (EMEA, MSCI EM) upcoming-event row.
# That record occupied row 2 of Sheet1 (directly under the header row);
# deleting it shifts every subsequent event row up by one, so the sheet
# goes from 16 data rows (A1:Q17) down to 15 data rows (A1:Q16) with no
# other content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:2").Delete()
